$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 146, col C: the "NA" placeholder text is removed, cell left blank
$ws.Range("C146").ClearContents()

# New rows 147-160: latest script run results for 2025-07-01
# Force column A to Text format first so the "YYYY-MM-DD" strings are not
# auto-converted into date serial numbers, matching the rest of the column.
$ws.Range("A147:A160").NumberFormat = "@"

$ws.Range("A147").Value = "2025-07-01"
$ws.Range("B147").Value = "développement durable"
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 1

$ws.Range("A148").Value = "2025-07-01"
$ws.Range("B148").Value = "zone tampon"
$ws.Range("C148").Value = 93
$ws.Range("D148").Value = 1

$ws.Range("A149").Value = "2025-07-01"
$ws.Range("B149").Value = "eaux souterraines"
$ws.Range("C149").Value = 93
$ws.Range("D149").Value = 1

$ws.Range("A150").Value = "2025-07-01"
$ws.Range("B150").Value = "ruissellement"
$ws.Range("C150").Value = 93
$ws.Range("D150").Value = 2

$ws.Range("A151").Value = "2025-07-01"
$ws.Range("B151").Value = "eaux de surface"
$ws.Range("C151").Value = 94
$ws.Range("D151").Value = 1

$ws.Range("A152").Value = "2025-07-01"
$ws.Range("B152").Value = "ruissellement"
$ws.Range("C152").Value = 94
$ws.Range("D152").Value = 1

$ws.Range("A153").Value = "2025-07-01"
$ws.Range("B153").Value = "ruissellement"
$ws.Range("C153").Value = 96
$ws.Range("D153").Value = 1

$ws.Range("A154").Value = "2025-07-01"
$ws.Range("B154").Value = "zone tampon"
$ws.Range("C154").Value = 98
$ws.Range("D154").Value = 5

$ws.Range("A155").Value = "2025-07-01"
$ws.Range("B155").Value = "eaux souterraines"
$ws.Range("C155").Value = 104
$ws.Range("D155").Value = 1

$ws.Range("A156").Value = "2025-07-01"
$ws.Range("B156").Value = "eaux souterraines"
$ws.Range("C156").Value = 105
$ws.Range("D156").Value = 3

$ws.Range("A157").Value = "2025-07-01"
$ws.Range("B157").Value = "eaux de surface"
$ws.Range("C157").Value = 106
$ws.Range("D157").Value = 4

$ws.Range("A158").Value = "2025-07-01"
$ws.Range("B158").Value = "eaux souterraines"
$ws.Range("C158").Value = 106
$ws.Range("D158").Value = 1

$ws.Range("A159").Value = "2025-07-01"
$ws.Range("B159").Value = "ruissellement"
$ws.Range("C159").Value = 110
$ws.Range("D159").Value = 1

$ws.Range("A160").Value = "2025-07-01"
$ws.Range("B160").Value = "développement durable"
$ws.Range("C160").Value = 113
$ws.Range("D160").Value = 1

# Reset the style on the new date cells back to the workbook default so no
# extra formatting is applied (only the underlying text format is kept).
$ws.Range("A147:A160").Style = "Normal"
